$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.780.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.287.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '102.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.856'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.290.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '43.673.76'
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '233.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("E23").Value = '  +14.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.58%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '177.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.41%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0359'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.42%  '
$ws.Range("E39").Value = '  -2.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.39%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("E45").Value = '  -4.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.102'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("E47").Value = '  +2.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.444'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.511.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.22%  '
